# Clean up the "Authors" column (column E) text: the raw author-list strings
# used a single space after each separating comma; the cleaned data adds one
# extra space after every comma. Apply this fix to every data row that has a
# non-empty author list (rows whose value is not "[]").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)   # column E = 5th column ("Authors")
    $old = $cell.Value()
    if ($old -ne $null -and $old -ne "[]") {
        $new = $old -replace ',(\s)', ', $1'
        if ($new -ne $old) {
            $cell.Value = $new
        }
    }
}
